# Boolean Peaking and Flexibility Flags.xlsx
# "update ets and bpaff to allow coal, nat gas to act as peakers"
#
# On the "BPaFF-BITPTaP" sheet (Is This Plant Type a Peaker), flag
# hard coal, natural gas (nonpeaker) and nuclear as peaker plant types
# (value 1 instead of 0). The "lignite" summary row (B13 = B2) follows
# automatically via its existing formula.
#
# Also bring the view state in line with the edited workbook: the
# "BPaFF-BITPTaP" sheet becomes the active/selected tab (with B5
# selected), and the "About" sheet's remembered selection moves to C4.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsPeaker = $wb.Worksheets.Item("BPaFF-BITPTaP")

# Update the boolean peaker flags.
$wsPeaker.Range("B2").Value = 1   # hard coal
$wsPeaker.Range("B3").Value = 1   # natural gas nonpeaker
$wsPeaker.Range("B4").Value = 1   # nuclear

# Restore/relocate the remembered cell selection on the "About" sheet.
$wsAbout.Activate() | Out-Null
$wsAbout.Range("C4").Select() | Out-Null

# Make "BPaFF-BITPTaP" the active sheet with B5 selected.
$wsPeaker.Activate() | Out-Null
$wsPeaker.Range("B5").Select() | Out-Null
